# Auto-generated edit script: applies numeric corrections to profit/price tables
# across multiple worksheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 106550
$ws.Range("J3").Value = 106550
$ws.Range("L3").Value = 106550
$ws.Range("N3").Value = -106778
$ws.Range("H29").Value = 2175.6365
$ws.Range("I29").Value = 1381.5555
$ws.Range("J29").Value = 5749
$ws.Range("K29").Value = 4144.666499999999
$ws.Range("L29").Value = 17247
$ws.Range("M29").Value = -3863.666499999999
$ws.Range("N29").Value = -17809
$ws.Range("H32").Value = 3953
$ws.Range("I32").Value = 3684
$ws.Range("K32").Value = 3684
$ws.Range("M32").Value = -3358
$ws.Range("H58").Value = 532.5
$ws.Range("J58").Value = 1200
$ws.Range("L58").Value = 3600
$ws.Range("N58").Value = -3900
$ws.Range("H62").Value = 17758.5
$ws.Range("I62").Value = 41037
$ws.Range("K62").Value = 41037
$ws.Range("M62").Value = -40413
$ws.Range("H65").Value = 17758.5
$ws.Range("I65").Value = 41037
$ws.Range("K65").Value = 205185
$ws.Range("M65").Value = -202065
$ws.Range("H74").Value = 999
$ws.Range("I74").Value = 999
$ws.Range("K74").Value = 999
$ws.Range("M74").Value = -63
$ws.Range("H77").Value = 999
$ws.Range("I77").Value = 999
$ws.Range("K77").Value = 4995
$ws.Range("M77").Value = -315
$ws.Range("H87").Value = 58749.5
$ws.Range("I87").Value = 47500
$ws.Range("J87").Value = 69999
$ws.Range("K87").Value = 47500
$ws.Range("L87").Value = 69999
$ws.Range("M87").Value = -46252
$ws.Range("N87").Value = -72495
$ws.Range("H90").Value = 58749.5
$ws.Range("I90").Value = 47500
$ws.Range("J90").Value = 69999
$ws.Range("K90").Value = 142500
$ws.Range("L90").Value = 209997
$ws.Range("M90").Value = -136260
$ws.Range("N90").Value = -222477
$ws.Range("H102").Value = 106550
$ws.Range("J102").Value = 106550
$ws.Range("L102").Value = 106550
$ws.Range("N102").Value = -113040
$ws.Range("H116").Value = 4981.6665
$ws.Range("I116").Value = 5065.7856
$ws.Range("J116").Value = 4813.4287
$ws.Range("K116").Value = 5065.7856
$ws.Range("L116").Value = 4813.4287
$ws.Range("M116").Value = -1623.7856
$ws.Range("N116").Value = -11697.4287
$ws.Range("H132").Value = 4327.683
$ws.Range("I132").Value = 3506.7837
$ws.Range("K132").Value = 10520.3511
$ws.Range("M132").Value = -7990.3511
$ws.Range("H137").Value = 5685.2856
$ws.Range("I137").Value = 1303
$ws.Range("K137").Value = 3909
$ws.Range("M137").Value = -1359
$ws.Range("H138").Value = 3055.2712
$ws.Range("J138").Value = 3083.7346
$ws.Range("L138").Value = 9251.203799999999
$ws.Range("N138").Value = -19531.2038

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 818.89746
$ws.Range("J2").Value = 1302.5
$ws.Range("L2").Value = 1302.5
$ws.Range("N2").Value = -1528.5
$ws.Range("H21").Value = 4000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H32").Value = 9829.932000000001
$ws.Range("I32").Value = 6044.6763
$ws.Range("J32").Value = 22699.8
$ws.Range("K32").Value = 6044.6763
$ws.Range("L32").Value = 22699.8
$ws.Range("M32").Value = -5757.6763
$ws.Range("N32").Value = -23273.8
$ws.Range("H45").Value = 2592.818
$ws.Range("I45").Value = 1185.75
$ws.Range("J45").Value = 3396.8572
$ws.Range("K45").Value = 1185.75
$ws.Range("L45").Value = 3396.8572
$ws.Range("M45").Value = -808.75
$ws.Range("N45").Value = -4150.8572
$ws.Range("H61").Value = 4030.9148
$ws.Range("I61").Value = 3718.525
$ws.Range("K61").Value = 3718.525
$ws.Range("M61").Value = -3506.525
$ws.Range("H74").Value = 126575.875
$ws.Range("I74").Value = 143515.28
$ws.Range("K74").Value = 143515.28
$ws.Range("M74").Value = -142641.28
$ws.Range("H77").Value = 126575.875
$ws.Range("I77").Value = 143515.28
$ws.Range("K77").Value = 717576.4
$ws.Range("M77").Value = -713208.4
$ws.Range("H88").Value = 2897.7
$ws.Range("I88").Value = 2885.5
$ws.Range("J88").Value = 2905.8333
$ws.Range("K88").Value = 2885.5
$ws.Range("L88").Value = 2905.8333
$ws.Range("M88").Value = -2479.5
$ws.Range("N88").Value = -3717.8333
$ws.Range("H91").Value = 2897.7
$ws.Range("I91").Value = 2885.5
$ws.Range("J91").Value = 2905.8333
$ws.Range("K91").Value = 2885.5
$ws.Range("L91").Value = 2905.8333
$ws.Range("M91").Value = -1481.5
$ws.Range("N91").Value = -5713.8333
$ws.Range("H109").Value = 64499
$ws.Range("J109").Value = 64499
$ws.Range("L109").Value = 64499
$ws.Range("N109").Value = -67273
$ws.Range("H116").Value = 818.89746
$ws.Range("J116").Value = 1302.5
$ws.Range("L116").Value = 1302.5
$ws.Range("N116").Value = -5890.5
$ws.Range("H124").Value = 64994
$ws.Range("J124").Value = 64994
$ws.Range("L124").Value = 64994
$ws.Range("N124").Value = -74814
$ws.Range("H131").Value = 86224.75
$ws.Range("J131").Value = 86224.75
$ws.Range("L131").Value = 86224.75
$ws.Range("N131").Value = -96304.75
$ws.Range("H132").Value = 34863.387
$ws.Range("I132").Value = 35758.832
$ws.Range("K132").Value = 107276.496
$ws.Range("M132").Value = -104746.496
$ws.Range("H136").Value = 4030.9148
$ws.Range("I136").Value = 3718.525
$ws.Range("K136").Value = 11155.575
$ws.Range("M136").Value = -8605.575000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 818.89746
$ws.Range("J3").Value = 1302.5
$ws.Range("L3").Value = 1302.5
$ws.Range("N3").Value = -1530.5
$ws.Range("H86").Value = 2371.9
$ws.Range("I86").Value = 2402.25
$ws.Range("J86").Value = 2351.6667
$ws.Range("K86").Value = 2402.25
$ws.Range("L86").Value = 2351.6667
$ws.Range("M86").Value = -1279.25
$ws.Range("N86").Value = -4597.6667
$ws.Range("H89").Value = 2371.9
$ws.Range("I89").Value = 2402.25
$ws.Range("J89").Value = 2351.6667
$ws.Range("K89").Value = 12011.25
$ws.Range("L89").Value = 11758.3335
$ws.Range("M89").Value = -6395.25
$ws.Range("N89").Value = -22990.3335
$ws.Range("H132").Value = 110000
$ws.Range("J132").Value = 110000
$ws.Range("L132").Value = 110000
$ws.Range("N132").Value = -120120
$ws.Range("H134").Value = 2101.7192
$ws.Range("I134").Value = 1869.7084
$ws.Range("K134").Value = 5609.1252
$ws.Range("M134").Value = -3074.1252

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3009.65
$ws.Range("I31").Value = 3009.65
$ws.Range("K31").Value = 3009.65
$ws.Range("M31").Value = -2714.65
$ws.Range("H34").Value = 3009.65
$ws.Range("I34").Value = 3009.65
$ws.Range("K34").Value = 3009.65
$ws.Range("M34").Value = -2807.65
$ws.Range("H59").Value = 93992.5
$ws.Range("J59").Value = 94990
$ws.Range("L59").Value = 94990
$ws.Range("N59").Value = -97280
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H97").Value = 52031.668
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 52031.668
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 52031.668
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -54013.668
$ws.Range("H99").Value = 2781.5789
$ws.Range("I99").Value = 1993
$ws.Range("J99").Value = 3865.875
$ws.Range("K99").Value = 1993
$ws.Range("L99").Value = 3865.875
$ws.Range("M99").Value = -495
$ws.Range("N99").Value = -6861.875
$ws.Range("H105").Value = 3147.5386
$ws.Range("I105").Value = 1927.625
$ws.Range("K105").Value = 1927.625
$ws.Range("M105").Value = -180.625
$ws.Range("H126").Value = 2781.5789
$ws.Range("I126").Value = 1993
$ws.Range("J126").Value = 3865.875
$ws.Range("K126").Value = 5979
$ws.Range("L126").Value = 11597.625
$ws.Range("M126").Value = -3509
$ws.Range("N126").Value = -16537.625
$ws.Range("H134").Value = 35072.773
$ws.Range("I134").Value = 39807.63
$ws.Range("J134").Value = 3112.5
$ws.Range("K134").Value = 119422.89
$ws.Range("L134").Value = 9337.5
$ws.Range("M134").Value = -116887.89
$ws.Range("N134").Value = -14407.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1183185.1
$ws.Range("I4").Value = 793667.9
$ws.Range("K4").Value = 2381003.7
$ws.Range("M4").Value = -2380891.7
$ws.Range("H29").Value = 1562.8049
$ws.Range("I29").Value = 1438.1666
$ws.Range("J29").Value = 1660.3478
$ws.Range("K29").Value = 4314.4998
$ws.Range("L29").Value = 4981.0434
$ws.Range("M29").Value = -4037.4998
$ws.Range("N29").Value = -5535.0434
$ws.Range("H117").Value = 3109
$ws.Range("J117").Value = 4328.4287
$ws.Range("L117").Value = 12985.2861
$ws.Range("N117").Value = -19869.2861
$ws.Range("H123").Value = 2198.5715
$ws.Range("I123").Value = 2198.5715
$ws.Range("K123").Value = 6595.7145
$ws.Range("M123").Value = -4145.7145
$ws.Range("H128").Value = 171919.42
$ws.Range("I128").Value = 171919.42
$ws.Range("K128").Value = 515758.26
$ws.Range("M128").Value = -510778.26
$ws.Range("H129").Value = 608795.5
$ws.Range("I129").Value = 1289.3334
$ws.Range("K129").Value = 3868.0002
$ws.Range("M129").Value = 1131.9998
$ws.Range("H131").Value = 4175781.8
$ws.Range("J131").Value = 6680265
$ws.Range("L131").Value = 20040795
$ws.Range("N131").Value = -20050875
$ws.Range("H133").Value = 8250
$ws.Range("I133").Value = 1000
$ws.Range("J133").Value = 10062.5
$ws.Range("K133").Value = 3000
$ws.Range("L133").Value = 30187.5
$ws.Range("M133").Value = 2060
$ws.Range("N133").Value = -40307.5
$ws.Range("H136").Value = 3597.8572
$ws.Range("I136").Value = 2692
$ws.Range("K136").Value = 8076
$ws.Range("M136").Value = -2976
$ws.Range("H138").Value = 1848.4
$ws.Range("J138").Value = 2500
$ws.Range("L138").Value = 7500
$ws.Range("N138").Value = -17780

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 121.05556
$ws.Range("I2").Value = 137
$ws.Range("J2").Value = 79.59999999999999
$ws.Range("K2").Value = 137
$ws.Range("L2").Value = 79.59999999999999
$ws.Range("M2").Value = -24
$ws.Range("N2").Value = -305.6
$ws.Range("H57").Value = 10500.6
$ws.Range("I57").Value = 8334.111000000001
$ws.Range("J57").Value = 29999
$ws.Range("K57").Value = 8334.111000000001
$ws.Range("L57").Value = 29999
$ws.Range("M57").Value = -7514.111000000001
$ws.Range("N57").Value = -31639
$ws.Range("H122").Value = 2133.76
$ws.Range("I122").Value = 1299.1765
$ws.Range("J122").Value = 3907.25
$ws.Range("K122").Value = 3897.5295
$ws.Range("L122").Value = 11721.75
$ws.Range("M122").Value = -1447.5295
$ws.Range("N122").Value = -16621.75
$ws.Range("H126").Value = 3783.5557
$ws.Range("I126").Value = 3175.279
$ws.Range("K126").Value = 9525.837
$ws.Range("M126").Value = -7055.837
$ws.Range("H132").Value = 25233.205
$ws.Range("I132").Value = 27465.691
$ws.Range("J132").Value = 7819.8
$ws.Range("K132").Value = 82397.073
$ws.Range("L132").Value = 23459.4
$ws.Range("M132").Value = -79867.073
$ws.Range("N132").Value = -28519.4
$ws.Range("H139").Value = 73999.25
$ws.Range("J139").Value = 73999.25
$ws.Range("L139").Value = 73999.25
$ws.Range("N139").Value = -84279.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13481.2
$ws.Range("I7").Value = 19138.666
$ws.Range("J7").Value = 4995
$ws.Range("K7").Value = 19138.666
$ws.Range("L7").Value = 4995
$ws.Range("M7").Value = -19026.666
$ws.Range("N7").Value = -5219
$ws.Range("H16").Value = 979.63635
$ws.Range("I16").Value = 786.7778
$ws.Range("K16").Value = 786.7778
$ws.Range("M16").Value = -616.7778
$ws.Range("H46").Value = 17403.95
$ws.Range("I46").Value = 20652.688
$ws.Range("J46").Value = 4409
$ws.Range("K46").Value = 20652.688
$ws.Range("L46").Value = 4409
$ws.Range("M46").Value = -20464.688
$ws.Range("N46").Value = -4785
$ws.Range("H55").Value = 1498
$ws.Range("J55").Value = 1842.875
$ws.Range("L55").Value = 1842.875
$ws.Range("N55").Value = -2188.875
$ws.Range("H68").Value = 3126.2
$ws.Range("I68").Value = 1943.6666
$ws.Range("K68").Value = 1943.6666
$ws.Range("M68").Value = -1194.6666
$ws.Range("H71").Value = 3126.2
$ws.Range("I71").Value = 1943.6666
$ws.Range("K71").Value = 9718.333000000001
$ws.Range("M71").Value = -5974.333000000001
$ws.Range("H93").Value = 1321.6428
$ws.Range("I93").Value = 1154.0769
$ws.Range("K93").Value = 1154.0769
$ws.Range("M93").Value = 93.92309999999998
$ws.Range("H100").Value = 4270.25
$ws.Range("I100").Value = 4126.8
$ws.Range("J100").Value = 4987.5
$ws.Range("K100").Value = 4126.8
$ws.Range("L100").Value = 4987.5
$ws.Range("M100").Value = -3585.8
$ws.Range("N100").Value = -6069.5
$ws.Range("H122").Value = 3708.4792
$ws.Range("J122").Value = 4335.5415
$ws.Range("L122").Value = 13006.6245
$ws.Range("N122").Value = -17906.6245
$ws.Range("H126").Value = 13481.2
$ws.Range("I126").Value = 19138.666
$ws.Range("J126").Value = 4995
$ws.Range("K126").Value = 57415.99800000001
$ws.Range("L126").Value = 14985
$ws.Range("M126").Value = -54945.99800000001
$ws.Range("N126").Value = -19925
$ws.Range("H127").Value = 121254
$ws.Range("J127").Value = 121254
$ws.Range("L127").Value = 121254
$ws.Range("N127").Value = -131174
$ws.Range("H132").Value = 33027.4
$ws.Range("I132").Value = 40049.875
$ws.Range("J132").Value = 4937.5
$ws.Range("K132").Value = 120149.625
$ws.Range("L132").Value = 14812.5
$ws.Range("M132").Value = -117619.625
$ws.Range("N132").Value = -19872.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 11033.167
$ws.Range("I4").Value = 20566.334
$ws.Range("K4").Value = 20566.334
$ws.Range("M4").Value = -20453.334
$ws.Range("H62").Value = 61935.5
$ws.Range("I62").Value = 3927
$ws.Range("J62").Value = 90939.75
$ws.Range("K62").Value = 3927
$ws.Range("L62").Value = 90939.75
$ws.Range("M62").Value = -3303
$ws.Range("N62").Value = -92187.75
$ws.Range("H65").Value = 61935.5
$ws.Range("I65").Value = 3927
$ws.Range("J65").Value = 90939.75
$ws.Range("K65").Value = 19635
$ws.Range("L65").Value = 454698.75
$ws.Range("M65").Value = -16515
$ws.Range("N65").Value = -460938.75
$ws.Range("H113").Value = 927.1429000000001
$ws.Range("I113").Value = 456.70587
$ws.Range("K113").Value = 1370.11761
$ws.Range("M113").Value = 799.88239
$ws.Range("H126").Value = 47344.684
$ws.Range("I126").Value = 85381.586
$ws.Range("J126").Value = 1700.4
$ws.Range("K126").Value = 256144.758
$ws.Range("L126").Value = 5101.200000000001
$ws.Range("M126").Value = -253674.758
$ws.Range("N126").Value = -10041.2
$ws.Range("H132").Value = 22669.771
$ws.Range("I132").Value = 24017.955
$ws.Range("K132").Value = 72053.86500000001
$ws.Range("M132").Value = -69523.86500000001
